# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (hunk 0)
$ws.Range("H28").Value = 1100.8334
$ws.Range("I28").Value = 1113.4375
$ws.Range("K28").Value = 1113.4375
$ws.Range("M28").Value = -628.4375
# Row 62 (hunk 1)
$ws.Range("H62").Value = 7035.1177
$ws.Range("J62").Value = 7120.1665
$ws.Range("L62").Value = 7120.1665
$ws.Range("N62").Value = -8368.166499999999
# Row 65 (hunk 2)
$ws.Range("H65").Value = 7035.1177
$ws.Range("J65").Value = 7120.1665
$ws.Range("L65").Value = 35600.8325
$ws.Range("N65").Value = -41840.8325
# Row 132 (hunk 3)
$ws.Range("H132").Value = 21283966
$ws.Range("I132").Value = 26318262
$ws.Range("K132").Value = 78954786
$ws.Range("M132").Value = -78952256
# Row 135 (hunk 4)
$ws.Range("H135").Value = 1665.44
$ws.Range("I135").Value = 501.7
$ws.Range("K135").Value = 4515.3
$ws.Range("M135").Value = -1980.3

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 5)
$ws.Range("H2").Value = 1092.2424
$ws.Range("I2").Value = 1087.56
$ws.Range("J2").Value = 1106.875
$ws.Range("K2").Value = 1087.56
$ws.Range("L2").Value = 1106.875
$ws.Range("M2").Value = -974.5599999999999
$ws.Range("N2").Value = -1332.875
# Row 32 (hunk 6)
$ws.Range("H32").Value = 8853.634
$ws.Range("I32").Value = 9402.916999999999
$ws.Range("K32").Value = 9402.916999999999
$ws.Range("M32").Value = -9115.916999999999
# Row 110 (hunk 7)
$ws.Range("H110").Value = 2020.2142
$ws.Range("I110").Value = 2029.4615
$ws.Range("K110").Value = 2029.4615
$ws.Range("M110").Value = 15.53850000000011
# Row 116 (hunk 8)
$ws.Range("H116").Value = 1092.2424
$ws.Range("I116").Value = 1087.56
$ws.Range("J116").Value = 1106.875
$ws.Range("K116").Value = 1087.56
$ws.Range("L116").Value = 1106.875
$ws.Range("M116").Value = 1206.44
$ws.Range("N116").Value = -5694.875
# Row 132 (hunk 9)
$ws.Range("H132").Value = 2297.3708
$ws.Range("I132").Value = 2146.4182
$ws.Range("K132").Value = 6439.2546
$ws.Range("M132").Value = -3909.2546

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 10)
$ws.Range("H3").Value = 1092.2424
$ws.Range("I3").Value = 1087.56
$ws.Range("J3").Value = 1106.875
$ws.Range("K3").Value = 1087.56
$ws.Range("L3").Value = 1106.875
$ws.Range("M3").Value = -973.5599999999999
$ws.Range("N3").Value = -1334.875
# Row 29 (hunk 11)
$ws.Range("H29").Value = 13047.25
$ws.Range("J29").Value = 12996
$ws.Range("L29").Value = 12996
$ws.Range("N29").Value = -13574
# Row 99 (hunk 12)
$ws.Range("H99").Value = 3008.6428
$ws.Range("I99").Value = 2552.6875
$ws.Range("J99").Value = 3616.5833
$ws.Range("K99").Value = 2552.6875
$ws.Range("L99").Value = 3616.5833
$ws.Range("M99").Value = -1054.6875
$ws.Range("N99").Value = -6612.5833
# Row 105 (hunk 13)
$ws.Range("H105").Value = 4496.737
$ws.Range("I105").Value = 2974.1428
$ws.Range("K105").Value = 2974.1428
$ws.Range("M105").Value = -1227.1428

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (hunk 14)
$ws.Range("H22").Value = 649.8889
$ws.Range("I22").Value = 645.2727
$ws.Range("J22").Value = 657.1429000000001
$ws.Range("K22").Value = 645.2727
$ws.Range("L22").Value = 657.1429000000001
$ws.Range("M22").Value = -295.2727
$ws.Range("N22").Value = -1357.1429
# Row 31 (hunk 15)
$ws.Range("H31").Value = 1928.7833
$ws.Range("I31").Value = 1135
$ws.Range("J31").Value = 3641.6843
$ws.Range("K31").Value = 1135
$ws.Range("L31").Value = 3641.6843
$ws.Range("M31").Value = -840
$ws.Range("N31").Value = -4231.6843
# Row 34 (hunk 16)
$ws.Range("H34").Value = 1928.7833
$ws.Range("I34").Value = 1135
$ws.Range("J34").Value = 3641.6843
$ws.Range("K34").Value = 1135
$ws.Range("L34").Value = 3641.6843
$ws.Range("M34").Value = -933
$ws.Range("N34").Value = -4045.6843
# Row 58 (hunk 17)
$ws.Range("H58").Value = 1954.0588
$ws.Range("I58").Value = 1936.0834
$ws.Range("J58").Value = 1997.2
$ws.Range("K58").Value = 1936.0834
$ws.Range("L58").Value = 1997.2
$ws.Range("M58").Value = -1733.0834
$ws.Range("N58").Value = -2403.2
# Row 96 (hunk 18)
$ws.Range("H96").Value = 10966
$ws.Range("J96").Value = 10966
$ws.Range("L96").Value = 10966
$ws.Range("N96").Value = -16458
# Row 132 (hunk 19)
$ws.Range("H132").Value = 1483938.9
$ws.Range("I132").Value = 1669119.6
$ws.Range("K132").Value = 5007358.800000001
$ws.Range("M132").Value = -5004828.800000001
# Row 134 (hunk 20)
$ws.Range("H134").Value = 2896.9836
$ws.Range("I134").Value = 1179.4884
$ws.Range("K134").Value = 3538.4652
$ws.Range("M134").Value = -1003.4652
# Row 136 (hunk 21)
$ws.Range("H136").Value = 1954.0588
$ws.Range("I136").Value = 1936.0834
$ws.Range("J136").Value = 1997.2
$ws.Range("K136").Value = 5808.2502
$ws.Range("L136").Value = 5991.6
$ws.Range("M136").Value = -3258.2502
$ws.Range("N136").Value = -11091.6

$ws = $wb.Worksheets.Item("CUL")
# Row 124 (hunk 22)
$ws.Range("H124").Value = 3563.8
$ws.Range("I124").Value = 3563.8
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 10691.4
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -5781.400000000001
$ws.Range("N124").ClearContents()
# Row 125 (hunk 23)
$ws.Range("H125").Value = 4508.875
$ws.Range("I125").Value = 3390.75
$ws.Range("J125").Value = 5627
$ws.Range("K125").Value = 10172.25
$ws.Range("L125").Value = 16881
$ws.Range("M125").Value = -5252.25
$ws.Range("N125").Value = -26721
# Row 126 (hunk 24)
$ws.Range("H126").Value = 4643.6665
$ws.Range("I126").Value = 1949
$ws.Range("J126").Value = 10033
$ws.Range("K126").Value = 5847
$ws.Range("L126").Value = 30099
$ws.Range("M126").Value = -907
$ws.Range("N126").Value = -39979

$ws = $wb.Worksheets.Item("GSM")
# Row 26 (hunk 25)
$ws.Range("H26").Value = 35555
$ws.Range("J26").Value = 35555
$ws.Range("L26").Value = 35555
$ws.Range("N26").Value = -36115
# Row 50 (hunk 26)
$ws.Range("H50").Value = 35555
$ws.Range("J50").Value = 35555
$ws.Range("L50").Value = 35555
$ws.Range("N50").Value = -36551
# Row 102 (hunk 27)
$ws.Range("H102").Value = 28893.945
$ws.Range("I102").Value = 1515.5834
$ws.Range("J102").Value = 79438.62
$ws.Range("K102").Value = 1515.5834
$ws.Range("L102").Value = 79438.62
$ws.Range("M102").Value = 106.4166
$ws.Range("N102").Value = -82682.62
# Row 113 (hunk 28)
$ws.Range("H113").Value = 1974.6666
$ws.Range("I113").Value = 1974.6666
$ws.Range("K113").Value = 1974.6666
$ws.Range("M113").Value = 195.3334
# Row 120 (hunk 29)
$ws.Range("H120").Value = 89997.5
$ws.Range("J120").Value = 89997.5
$ws.Range("L120").Value = 89997.5
$ws.Range("N120").Value = -99673.5

$ws = $wb.Worksheets.Item("LTW")
# Row 14 (hunk 30)
$ws.Range("H14").Value = 9999.333000000001
$ws.Range("J14").Value = 9999.333000000001
$ws.Range("L14").Value = 9999.333000000001
$ws.Range("N14").Value = -10343.333
# Row 22 (hunk 31)
$ws.Range("H22").Value = 4110.8887
$ws.Range("J22").Value = 4110.8887
$ws.Range("L22").Value = 4110.8887
$ws.Range("N22").Value = -4700.8887
# Row 27 (hunk 32)
$ws.Range("H27").Value = 4110.8887
$ws.Range("J27").Value = 4110.8887
$ws.Range("L27").Value = 4110.8887
$ws.Range("N27").Value = -4324.8887
# Row 46 (hunk 33)
$ws.Range("H46").Value = 2623.652
$ws.Range("J46").Value = 2596
$ws.Range("L46").Value = 2596
$ws.Range("N46").Value = -2972
# Row 55 (hunk 34)
$ws.Range("H55").Value = 253.17647
$ws.Range("I55").Value = 259.7
$ws.Range("K55").Value = 259.7
$ws.Range("M55").Value = -86.69999999999999
# Row 61 (hunk 35)
$ws.Range("H61").Value = 4485.722
$ws.Range("I61").Value = 4677.6875
$ws.Range("K61").Value = 4677.6875
$ws.Range("M61").Value = -4475.6875
# Row 113 (hunk 36)
$ws.Range("H113").Value = 4485.722
$ws.Range("I113").Value = 4677.6875
$ws.Range("K113").Value = 4677.6875
$ws.Range("M113").Value = -2507.6875
# Row 121 (hunk 37)
$ws.Range("H121").Value = 94998.5
$ws.Range("J121").Value = 89998
$ws.Range("L121").Value = 89998
$ws.Range("N121").Value = -93492
# Row 122 (hunk 38)
$ws.Range("H122").Value = 2890.7827
$ws.Range("I122").Value = 2174.5
$ws.Range("K122").Value = 6523.5
$ws.Range("M122").Value = -4073.5

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (hunk 39)
$ws.Range("H96").Value = 3749
$ws.Range("I96").Value = 2446.5
$ws.Range("K96").Value = 2446.5
$ws.Range("M96").Value = -1073.5
# Row 100 (hunk 40)
$ws.Range("H100").Value = 3089.3333
$ws.Range("I100").Value = 2465.75
$ws.Range("J100").Value = 4336.5
$ws.Range("K100").Value = 4931.5
$ws.Range("L100").Value = 8673
$ws.Range("M100").Value = -4390.5
$ws.Range("N100").Value = -9755
# Row 113 (hunk 41)
$ws.Range("H113").Value = 3473110
$ws.Range("I113").Value = 6410745.5
$ws.Range("J113").Value = 1359.1818
$ws.Range("K113").Value = 19232236.5
$ws.Range("L113").Value = 4077.5454
$ws.Range("M113").Value = -19230066.5
$ws.Range("N113").Value = -8417.545399999999

